$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.158808827400208
$ws.Range("B1").Value = 2.748419523239136
$ws.Range("C1").Value = 6.942601203918457
$ws.Range("D1").Value = 1.988441944122314
$ws.Range("E1").Value = 1.113291501998901
